$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row8 = New-Object 'object[,]' 1,29
$row8[0,0] = 6
$row8[0,1] = 7095169
$row8[0,2] = "Greece Super League 1"
$row8[0,3] = "Greece Super League 1"
$row8[0,4] = 45164.58333333334
$row8[0,5] = "Asteras Tripolis"
$row8[0,6] = "OFI Crete"
$row8[0,7] = 3
$row8[0,8] = 0
$row8[0,9] = "H"
$row8[0,10] = 2.375
$row8[0,11] = 3.1
$row8[0,12] = 3.1
$row8[0,13] = 2.3
$row8[0,14] = 3.1
$row8[0,15] = 3.4
$row8[0,16] = -0.25
$row8[0,17] = 1.95
$row8[0,18] = 1.9
$row8[0,19] = 2
$row8[0,20] = 1.875
$row8[0,21] = 1.975
$row8[0,22] = 1.3
$row8[0,23] = -1
$row8[0,24] = -1
$row8[0,25] = 0.95
$row8[0,26] = -1
$row8[0,27] = 0.875
$row8[0,28] = -1
$ws.Range("A8:AC8").Value2 = $row8

$row9 = New-Object 'object[,]' 1,29
$row9[0,0] = 7
$row9[0,1] = 7095328
$row9[0,2] = "Greece Super League 1"
$row9[0,3] = "Greece Super League 1"
$row9[0,4] = 45164.58333333334
$row9[0,5] = "Panathinaikos"
$row9[0,6] = "Volos NFC"
$row9[0,7] = 3
$row9[0,8] = 0
$row9[0,9] = "H"
$row9[0,10] = 1.3
$row9[0,11] = 5
$row9[0,12] = 11
$row9[0,13] = 1.25
$row9[0,14] = 5.5
$row9[0,15] = 13
$row9[0,16] = -1.75
$row9[0,17] = 2
$row9[0,18] = 1.85
$row9[0,19] = 2.75
$row9[0,20] = 1.925
$row9[0,21] = 1.925
$row9[0,22] = 0.25
$row9[0,23] = -1
$row9[0,24] = -1
$row9[0,25] = 1
$row9[0,26] = -1
$row9[0,27] = 0.4625
$row9[0,28] = -0.5
$ws.Range("A9:AC9").Value2 = $row9

$row124 = New-Object 'object[,]' 1,29
$row124[0,0] = 122
$row124[0,1] = 6937238
$row124[0,2] = "Greece Super League 1"
$row124[0,3] = "Greece Super League 1"
$row124[0,4] = 45305.64583333334
$row124[0,5] = "PAOK Salonika"
$row124[0,6] = "Giannina"
$row124[0,7] = 4
$row124[0,8] = 0
$row124[0,9] = "H"
$row124[0,10] = 1.111
$row124[0,11] = 9
$row124[0,12] = 23
$row124[0,13] = 1.25
$row124[0,14] = 6
$row124[0,15] = 9
$row124[0,16] = -1.75
$row124[0,17] = 2.025
$row124[0,18] = 1.825
$row124[0,19] = 2.75
$row124[0,20] = 1.8
$row124[0,21] = 2.05
$row124[0,22] = 0.25
$row124[0,23] = -1
$row124[0,24] = -1
$row124[0,25] = 1.025
$row124[0,26] = -1
$row124[0,27] = 0.8
$row124[0,28] = -1
$ws.Range("A124:AC124").Value2 = $row124

$row125 = New-Object 'object[,]' 1,29
$row125[0,0] = 123
$row125[0,1] = 6936857
$row125[0,2] = "Greece Super League 1"
$row125[0,3] = "Greece Super League 1"
$row125[0,4] = 45305.64583333334
$row125[0,5] = "AEK Athens"
$row125[0,6] = "Panathinaikos"
$row125[0,7] = 2
$row125[0,8] = 2
$row125[0,9] = "D"
$row125[0,10] = 1.909
$row125[0,11] = 3.5
$row125[0,12] = 4.2
$row125[0,13] = 2.15
$row125[0,14] = 3.2
$row125[0,15] = 3.5
$row125[0,16] = -0.25
$row125[0,17] = 1.85
$row125[0,18] = 2
$row125[0,19] = 2
$row125[0,20] = 1.8
$row125[0,21] = 2.05
$row125[0,22] = -1
$row125[0,23] = 2.2
$row125[0,24] = -1
$row125[0,25] = -0.5
$row125[0,26] = 0.5
$row125[0,27] = 0.8
$row125[0,28] = -1
$ws.Range("A125:AC125").Value2 = $row125

$row168 = New-Object 'object[,]' 1,29
$row168[0,0] = 166
$row168[0,1] = 6937267
$row168[0,2] = "Greece Super League 1"
$row168[0,3] = "Greece Super League 1"
$row168[0,4] = 45350.41666666666
$row168[0,5] = "Volos NFC"
$row168[0,6] = "OFI Crete"
$row168[0,7] = 3
$row168[0,8] = 1
$row168[0,9] = "H"
$row168[0,10] = 2.7
$row168[0,11] = 3.25
$row168[0,12] = 2.625
$row168[0,13] = 2.7
$row168[0,14] = 3.2
$row168[0,15] = 2.8
$row168[0,16] = 0
$row168[0,17] = 1.825
$row168[0,18] = 2.025
$row168[0,19] = 2.25
$row168[0,20] = 2
$row168[0,21] = 1.85
$row168[0,22] = 1.7
$row168[0,23] = -1
$row168[0,24] = -1
$row168[0,25] = 0.825
$row168[0,26] = -1
$row168[0,27] = 1
$row168[0,28] = -1
$ws.Range("A168:AC168").Value2 = $row168

$row169 = New-Object 'object[,]' 1,29
$row169[0,0] = 167
$row169[0,1] = 6935703
$row169[0,2] = "Greece Super League 1"
$row169[0,3] = "Greece Super League 1"
$row169[0,4] = 45350.41666666666
$row169[0,5] = "Asteras Tripolis"
$row169[0,6] = "Kifisias FC"
$row169[0,7] = 3
$row169[0,8] = 3
$row169[0,9] = "D"
$row169[0,10] = 1.833
$row169[0,11] = 3.4
$row169[0,12] = 4.5
$row169[0,13] = 1.8
$row169[0,14] = 3.5
$row169[0,15] = 4.75
$row169[0,16] = -0.75
$row169[0,17] = 2.025
$row169[0,18] = 1.825
$row169[0,19] = 2.5
$row169[0,20] = 1.875
$row169[0,21] = 1.975
$row169[0,22] = -1
$row169[0,23] = 2.5
$row169[0,24] = -1
$row169[0,25] = -1
$row169[0,26] = 0.825
$row169[0,27] = 0.875
$row169[0,28] = -1
$ws.Range("A169:AC169").Value2 = $row169

$row189 = New-Object 'object[,]' 1,29
$row189[0,0] = 187
$row189[0,1] = 7920472
$row189[0,2] = "Greece Super League 1"
$row189[0,3] = "Greece Super League 1"
$row189[0,4] = 45382.5
$row189[0,5] = "Panathinaikos"
$row189[0,6] = "PAOK Salonika"
$row189[0,7] = $null
$row189[0,8] = $null
$row189[0,9] = $null
$row189[0,10] = 2.25
$row189[0,11] = 3.3
$row189[0,12] = 3.3
$row189[0,13] = 2.1
$row189[0,14] = 3.3
$row189[0,15] = 3.6
$row189[0,16] = -0.25
$row189[0,17] = 1.85
$row189[0,18] = 2
$row189[0,19] = 2.5
$row189[0,20] = 2.05
$row189[0,21] = 1.8
$row189[0,22] = 0
$row189[0,23] = 0
$row189[0,24] = 0
$row189[0,25] = 0
$row189[0,26] = 0
$row189[0,27] = $null
$row189[0,28] = $null
$ws.Range("A189:AC189").Value2 = $row189

$row190 = New-Object 'object[,]' 1,29
$row190[0,0] = 188
$row190[0,1] = 7920471
$row190[0,2] = "Greece Super League 1"
$row190[0,3] = "Greece Super League 1"
$row190[0,4] = 45382.60416666666
$row190[0,5] = "Aris Salonika"
$row190[0,6] = "Lamia"
$row190[0,7] = $null
$row190[0,8] = $null
$row190[0,9] = $null
$row190[0,10] = 1.571
$row190[0,11] = 4
$row190[0,12] = 6
$row190[0,13] = 1.5
$row190[0,14] = 4
$row190[0,15] = 7.5
$row190[0,16] = -1
$row190[0,17] = 1.85
$row190[0,18] = 2
$row190[0,19] = 2.5
$row190[0,20] = 2
$row190[0,21] = 1.85
$row190[0,22] = 0
$row190[0,23] = 0
$row190[0,24] = 0
$row190[0,25] = 0
$row190[0,26] = 0
$row190[0,27] = $null
$row190[0,28] = $null
$ws.Range("A190:AC190").Value2 = $row190

$row191 = New-Object 'object[,]' 1,29
$row191[0,0] = 189
$row191[0,1] = 7920470
$row191[0,2] = "Greece Super League 1"
$row191[0,3] = "Greece Super League 1"
$row191[0,4] = 45382.60416666666
$row191[0,5] = "AEK Athens"
$row191[0,6] = "Olympiakos"
$row191[0,7] = $null
$row191[0,8] = $null
$row191[0,9] = $null
$row191[0,10] = 1.909
$row191[0,11] = 3.4
$row191[0,12] = 4.2
$row191[0,13] = 1.8
$row191[0,14] = 3.4
$row191[0,15] = 5
$row191[0,16] = -0.5
$row191[0,17] = 1.8
$row191[0,18] = 2.05
$row191[0,19] = 2.5
$row191[0,20] = 1.975
$row191[0,21] = 1.875
$row191[0,22] = 0
$row191[0,23] = 0
$row191[0,24] = 0
$row191[0,25] = 0
$row191[0,26] = 0
$row191[0,27] = $null
$row191[0,28] = $null
$ws.Range("A191:AC191").Value2 = $row191

# Remove the trailing 4 rows that no longer exist in the updated feed (192-195)
$ws.Range("A192:A195").EntireRow.Delete()
